$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.23455854456145
$ws.Range("C2").Value = 13.87391521980242
$ws.Range("D2").Value = 15.08546375431524
$ws.Range("E2").Value = 16.51362543263631
$ws.Range("G2").Value = 54.23324736013105
$ws.Range("H2").Value = 20.1162518101791
$ws.Range("J2").Value = 9.43852394366068
$ws.Range("N2").Value = 18.85098631877658
$ws.Range("B3").Value = 20.738566854177
$ws.Range("C3").Value = 13.41222186052857
$ws.Range("D3").Value = 15.03386908170059
$ws.Range("E3").Value = 16.46319788209662
$ws.Range("G3").Value = 53.78621963549934
$ws.Range("H3").Value = 20.11172354976033
$ws.Range("J3").Value = 9.451990614089588
$ws.Range("N3").Value = 18.92887399361645
$ws.Range("B4").Value = 20.43442100508006
$ws.Range("C4").Value = 13.12473643310034
$ws.Range("D4").Value = 15.00590092595192
$ws.Range("E4").Value = 16.43632105992922
$ws.Range("G4").Value = 53.52959402250321
$ws.Range("H4").Value = 20.11411857892905
$ws.Range("J4").Value = 9.461856880183317
$ws.Range("N4").Value = 18.9787319142674
$ws.Range("B5").Value = 20.3107739834865
$ws.Range("C5").Value = 13.00677840484208
$ws.Range("D5").Value = 14.99544252470356
$ws.Range("E5").Value = 16.42640083614445
$ws.Range("G5").Value = 53.42959092383244
$ws.Range("H5").Value = 20.11639202504955
$ws.Range("J5").Value = 9.466278609880071
$ws.Range("N5").Value = 18.99956271619718
$ws.Range("B6").Value = 20.29026610862549
$ws.Range("C6").Value = 12.98714909521884
$ws.Range("D6").Value = 14.99376276198316
$ws.Range("E6").Value = 16.42481605805858
$ws.Range("G6").Value = 53.41326401867401
$ws.Range("H6").Value = 20.11684776055416
$ws.Range("J6").Value = 9.467037044148862
$ws.Range("N6").Value = 19.00305270617094
$ws.Range("B7").Value = 20.43275200230628
$ws.Range("C7").Value = 13.12314859379587
$ws.Range("D7").Value = 15.00575607232103
$ws.Range("E7").Value = 16.43618308702205
$ws.Range("G7").Value = 53.52822672187931
$ws.Range("H7").Value = 20.11414399210051
$ws.Range("J7").Value = 9.461914889804589
$ws.Range("N7").Value = 18.97901076530376
$ws.Range("B8").Value = 21.06357428167022
$ws.Range("C8").Value = 13.7156705538815
$ws.Range("D8").Value = 15.06690737561417
$ws.Range("E8").Value = 16.49539325209834
$ws.Range("G8").Value = 54.0754613934708
$ws.Range("H8").Value = 20.11361417840964
$ws.Range("J8").Value = 9.442835307929794
$ws.Range("N8").Value = 18.87742073328032
$ws.Range("B9").Value = 22.29546180594446
$ws.Range("C9").Value = 14.83776141226997
$ws.Range("D9").Value = 15.21597397923338
$ws.Range("E9").Value = 16.64365910351539
$ws.Range("G9").Value = 55.2859999757682
$ws.Range("H9").Value = 20.15378189845859
$ws.Range("J9").Value = 9.418122632053327
$ws.Range("N9").Value = 18.69426771351791
$ws.Range("B10").Value = 23.18683822493992
$ws.Range("C10").Value = 15.62817411778227
$ws.Range("D10").Value = 15.34279515473439
$ws.Range("E10").Value = 16.77176555402118
$ws.Range("G10").Value = 56.25285974931048
$ws.Range("H10").Value = 20.20855831442798
$ws.Range("J10").Value = 9.407742603057539
$ws.Range("N10").Value = 18.5693859254463
$ws.Range("B11").Value = 23.58732951699473
$ws.Range("C11").Value = 15.97863860453613
$ws.Range("D11").Value = 15.40412475773099
$ws.Range("E11").Value = 16.83409421635852
$ws.Range("G11").Value = 56.70803242440576
$ws.Range("H11").Value = 20.23897636108053
$ws.Range("J11").Value = 9.404714739845652
$ws.Range("N11").Value = 18.51465192777622
$ws.Range("B12").Value = 23.73810999573579
$ws.Range("C12").Value = 16.10991793628135
$ws.Range("D12").Value = 15.42785996511346
$ws.Range("E12").Value = 16.85826783069855
$ws.Range("G12").Value = 56.88246589925103
$ws.Range("H12").Value = 20.25128518734216
$ws.Range("J12").Value = 9.403812090524356
$ws.Range("N12").Value = 18.49422225004045
$ws.Range("B13").Value = 23.70567811476815
$ws.Range("C13").Value = 16.08171020895953
$ws.Range("D13").Value = 15.4227256386028
$ws.Range("E13").Value = 16.85303640601835
$ws.Range("G13").Value = 56.84480871603451
$ws.Range("H13").Value = 20.2485991394634
$ws.Range("J13").Value = 9.403995638570372
$ws.Range("N13").Value = 18.49860896623488
$ws.Range("B14").Value = 23.59975275581607
$ws.Range("C14").Value = 15.9894684420616
$ws.Range("D14").Value = 15.40606729865429
$ws.Range("E14").Value = 16.8360716185199
$ws.Range("G14").Value = 56.72234249619881
$ws.Range("H14").Value = 20.23997318680621
$ws.Range("J14").Value = 9.404635587853843
$ws.Range("N14").Value = 18.51296522432886
$ws.Range("B15").Value = 23.5347517265546
$ws.Range("C15").Value = 15.93277746637547
$ws.Range("D15").Value = 15.39592975859642
$ws.Range("E15").Value = 16.82575420876421
$ws.Range("G15").Value = 56.64759376575781
$ws.Range("H15").Value = 20.23479240836155
$ws.Range("J15").Value = 9.405059351599864
$ws.Range("N15").Value = 18.52179746670681
$ws.Range("B16").Value = 23.16055080910059
$ws.Range("C16").Value = 15.60507629092488
$ws.Range("D16").Value = 15.33885931179033
$ws.Range("E16").Value = 16.76777280286548
$ws.Range("G16").Value = 56.22341029635308
$ws.Range("H16").Value = 20.20668105689753
$ws.Range("J16").Value = 9.407974600090521
$ws.Range("N16").Value = 18.57300455085229
$ws.Range("B17").Value = 22.92959411111328
$ws.Range("C17").Value = 15.40161797187871
$ws.Range("D17").Value = 15.30477189354601
$ws.Range("E17").Value = 16.73323305161226
$ws.Range("G17").Value = 55.96702140029788
$ws.Range("H17").Value = 20.19084411246996
$ws.Range("J17").Value = 9.410197161633493
$ws.Range("N17").Value = 18.60494888986634
$ws.Range("B18").Value = 22.79629201205108
$ws.Range("C18").Value = 15.28374432966195
$ws.Range("D18").Value = 15.28550892596899
$ws.Range("E18").Value = 16.71374889252521
$ws.Range("G18").Value = 55.82100600659445
$ws.Range("H18").Value = 20.18225290713328
$ws.Range("J18").Value = 9.411634958529397
$ws.Range("N18").Value = 18.62351786176859
$ws.Range("B19").Value = 22.75108387038272
$ws.Range("C19").Value = 15.24369245086453
$ws.Range("D19").Value = 15.27904611761543
$ws.Range("E19").Value = 16.70721787071958
$ws.Range("G19").Value = 55.77182130180025
$ws.Range("H19").Value = 20.17943300953782
$ws.Range("J19").Value = 9.412149144142976
$ws.Range("N19").Value = 18.62983860933525
$ws.Range("B20").Value = 22.95422883140641
$ws.Range("C20").Value = 15.42336536387056
$ws.Range("D20").Value = 15.30836512098638
$ws.Range("E20").Value = 16.73687039270511
$ws.Range("G20").Value = 55.99416506211987
$ws.Range("H20").Value = 20.19247639561878
$ws.Range("J20").Value = 9.40994406226134
$ws.Range("N20").Value = 18.60152814504309
$ws.Range("B21").Value = 23.63089060180914
$ws.Range("C21").Value = 16.01660193316576
$ws.Range("D21").Value = 15.41094648917049
$ws.Range("E21").Value = 16.84103918584261
$ws.Range("G21").Value = 56.75825874199904
$ws.Range("H21").Value = 20.24248540344456
$ws.Range("J21").Value = 9.404440996591697
$ws.Range("N21").Value = 18.50874039770817
$ws.Range("B22").Value = 24.06795181399221
$ws.Range("C22").Value = 16.39590965987262
$ws.Range("D22").Value = 15.48096169355625
$ws.Range("E22").Value = 16.91244226036465
$ws.Range("G22").Value = 57.26963415327441
$ws.Range("H22").Value = 20.27977392501411
$ws.Range("J22").Value = 9.402266375354563
$ws.Range("N22").Value = 18.44982809623495
$ws.Range("B23").Value = 23.83520688637662
$ws.Range("C23").Value = 16.19427307332319
$ws.Range("D23").Value = 15.44332552857665
$ws.Range("E23").Value = 16.87403315745418
$ws.Range("G23").Value = 56.99565215219707
$ws.Range("H23").Value = 20.2594514803755
$ws.Range("J23").Value = 9.403296812988811
$ws.Range("N23").Value = 18.48111292118791
$ws.Range("B24").Value = 22.94309308822886
$ws.Range("C24").Value = 15.41353617011176
$ws.Range("D24").Value = 15.30673958054783
$ws.Range("E24").Value = 16.73522478748025
$ws.Range("G24").Value = 55.98188908520751
$ws.Range("H24").Value = 20.19173684041183
$ws.Range("J24").Value = 9.410057990097332
$ws.Range("N24").Value = 18.6030740290385
$ws.Range("B25").Value = 21.96387922905013
$ws.Range("C25").Value = 14.53953445041705
$ws.Range("D25").Value = 15.17257095509634
$ws.Range("E25").Value = 16.60014642798474
$ws.Range("G25").Value = 54.94442351037438
$ws.Range("H25").Value = 20.13848557922777
$ws.Range("J25").Value = 9.423444544303239
$ws.Range("N25").Value = 18.74210732449649
